# Trading update: 2026-02-17 20:43:47
# Appends the new closed/open MarketMaking trade (Trade #63) to both the
# "All Trades" log and the per-strategy "MarketMaking" log.
#
# Helper: writes a cell as literal text (never let Excel's COM type
# inference reinterpret date-looking / numeric-looking strings), then
# resets the cell style back to Normal so no stray NumberFormat /
# quotePrefix style sticks around on the new row.
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = [char]39 + $text
    $cell.Style = "Normal"
}

function Set-NumberCell($ws, $row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "All Trades" sheet - new row 64
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$r = 64

Set-NumberCell $allTrades $r 1 63
Set-TextCell   $allTrades $r 2 "2026-02-17"
Set-TextCell   $allTrades $r 3 "20:42:37"
Set-TextCell   $allTrades $r 4 "MarketMaking"
Set-TextCell   $allTrades $r 5 "DOWN"
Set-NumberCell $allTrades $r 6 0.51
Set-TextCell   $allTrades $r 7 ""
Set-TextCell   $allTrades $r 8 "OPEN"
Set-NumberCell $allTrades $r 9 0
Set-NumberCell $allTrades $r 10 0
Set-NumberCell $allTrades $r 11 99.98999999999999
Set-TextCell   $allTrades $r 12 ""
Set-NumberCell $allTrades $r 13 0
Set-NumberCell $allTrades $r 14 0
Set-NumberCell $allTrades $r 15 0
Set-NumberCell $allTrades $r 16 0.6
Set-TextCell   $allTrades $r 17 "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------
# "MarketMaking" sheet - new row 31 (same trade, strategy-specific log
# uses a different column layout: L/M/N/O/P/Q instead of L..Q)
# ---------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$r2 = 31

Set-NumberCell $mm $r2 1 63
Set-TextCell   $mm $r2 2 "2026-02-17"
Set-TextCell   $mm $r2 3 "20:42:37"
Set-TextCell   $mm $r2 4 "MarketMaking"
Set-TextCell   $mm $r2 5 "DOWN"
Set-NumberCell $mm $r2 6 0.51
Set-TextCell   $mm $r2 7 ""
Set-TextCell   $mm $r2 8 "OPEN"
Set-NumberCell $mm $r2 9 0
Set-NumberCell $mm $r2 10 0
Set-NumberCell $mm $r2 11 99.98999999999999
Set-NumberCell $mm $r2 12 0
Set-NumberCell $mm $r2 13 0
Set-NumberCell $mm $r2 14 0.6
Set-TextCell   $mm $r2 15 "Normal spread capture: 19600 bps"
Set-TextCell   $mm $r2 16 ""
Set-NumberCell $mm $r2 17 0
